$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: updated rejection record ---
$ws.Range("B2").Formula = "'74399724366101545557798"
$ws.Range("D2").Value = "30/12/2024"
$ws.Range("E2").Value = 425115

# --- Row 3: clear the now-empty detail columns, keep Filiale/Motif ---
$ws.Range("B3:E3").ClearContents()

# --- Row 4: replace with the new "SG - COTE D IVOIRE" record ---
$ws.Range("A4").Value = "SG - COTE D IVOIRE"
$ws.Range("B4").Formula = "'74637414366101546447445"
$ws.Range("C4").Value = "TKNVHL"
$ws.Range("D4").Value = "26/12/2024"
$ws.Range("E4").Value = 92000
$ws.Range("F4").Value = "V0202 POS Entry Mode is ** invalid"

# --- Row 5: removed entirely ---
$ws.Range("A5:F5").ClearContents()

# --- Shrink the table / autofilter to the new data extent ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F4"))

# --- Column width tweaks (stored width = ColumnWidth + 0.8333333333333334) ---
$ws.Columns.Item(1).ColumnWidth = 19.166666666666668
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
